# Update the Navigation Smoke scripts testdata:
# Fix typo in the "LIVE" URL for the ipacs environment:
#   https://trail.provanaipacs.com/  ->  https://trial.provanaipacs.com/
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ipacs")
$ws.Range("B2").Value = "https://trial.provanaipacs.com/"
